# Generate Report for Handback
#
# This mirrors a localization "handback" run being recorded into the
# status workbook:
#   - the status text moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every file row (both language
#     sheets, which drives the same text on the Overview sheet as well,
#     since it reuses the same status string),
#   - the "Latest Handback DateTime" column (H) is stamped with the time
#     the handback finished (different per language sheet),
#   - new "Latest Target File" (F) / "Latest Handback File" (G) columns
#     are populated with hyperlinks pointing at the same file that was
#     handed off, recording that the file has now also been handed back.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: keep the summary rows showing the same status text ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# Helper hyperlink targets (same external file already referenced by the
# corresponding "Source"/"Latest Handoff File" cells in columns A and D).
$md911Target  = "https://github.com/OpenLocalizationTest/oltest/blob/53edfee334e575533e9fc207c8a74706ae780b34/e2e/911c866f-3374-4099-b8f3-5a6de249a572.md"
$md911Name    = "911c866f-3374-4099-b8f3-5a6de249a572.md"

$mdCaTarget   = "https://github.com/OpenLocalizationTest/oltest/blob/53edfee334e575533e9fc207c8a74706ae780b34/e2e/ca725fb2-9b93-441e-a073-4242e3867702.md"
$mdCaName     = "ca725fb2-9b93-441e-a073-4242e3867702.md"

$zh911Target  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c32ba1eb2f7d0461370aff6e31f3df5f10c22f73/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.zh-cn.xlf"
$zh911Name    = "911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.zh-cn.xlf"

$zhCaTarget   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c32ba1eb2f7d0461370aff6e31f3df5f10c22f73/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ca725fb2-9b93-441e-a073-4242e3867702.9957bbb691565fe1dd424a9939361ae51cc5e353.zh-cn.xlf"
$zhCaName     = "ca725fb2-9b93-441e-a073-4242e3867702.9957bbb691565fe1dd424a9939361ae51cc5e353.zh-cn.xlf"

$de911Target  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/021ddf4dbd035d5a0c00fea2851dc6227f5f8b84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.de-de.xlf"
$de911Name    = "911c866f-3374-4099-b8f3-5a6de249a572.cf28377ec91c04a903c845dbed1b68185b3cfd30.de-de.xlf"

$deCaTarget   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/021ddf4dbd035d5a0c00fea2851dc6227f5f8b84/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ca725fb2-9b93-441e-a073-4242e3867702.9957bbb691565fe1dd424a9939361ae51cc5e353.de-de.xlf"
$deCaName     = "ca725fb2-9b93-441e-a073-4242e3867702.9957bbb691565fe1dd424a9939361ae51cc5e353.de-de.xlf"

function Style-AsHyperlink($range) {
    # Match the workbook's existing custom "HyperLink" cell style
    # (underlined Calibri 11 in #6495ED) used by columns A/D.
    $range.Font.Underline = 2
    $range.Font.Color = 6591981
}

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("F2").Value = $md911Name
Style-AsHyperlink $zhcn.Range("F2")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $md911Target, "", "", $md911Name) | Out-Null

$zhcn.Range("G2").Value = $zh911Name
Style-AsHyperlink $zhcn.Range("G2")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), $zh911Target, "", "", $zh911Name) | Out-Null

$zhcn.Range("H2").Value = "2016-03-22 00:43:29"

$zhcn.Range("F3").Value = $mdCaName
Style-AsHyperlink $zhcn.Range("F3")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $mdCaTarget, "", "", $mdCaName) | Out-Null

$zhcn.Range("G3").Value = $zhCaName
Style-AsHyperlink $zhcn.Range("G3")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), $zhCaTarget, "", "", $zhCaName) | Out-Null

$zhcn.Range("H3").Value = "2016-03-22 00:43:29"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("F2").Value = $md911Name
Style-AsHyperlink $dede.Range("F2")
$dede.Hyperlinks.Add($dede.Range("F2"), $md911Target, "", "", $md911Name) | Out-Null

$dede.Range("G2").Value = $de911Name
Style-AsHyperlink $dede.Range("G2")
$dede.Hyperlinks.Add($dede.Range("G2"), $de911Target, "", "", $de911Name) | Out-Null

$dede.Range("H2").Value = "2016-03-22 00:43:36"

$dede.Range("F3").Value = $mdCaName
Style-AsHyperlink $dede.Range("F3")
$dede.Hyperlinks.Add($dede.Range("F3"), $mdCaTarget, "", "", $mdCaName) | Out-Null

$dede.Range("G3").Value = $deCaName
Style-AsHyperlink $dede.Range("G3")
$dede.Hyperlinks.Add($dede.Range("G3"), $deCaTarget, "", "", $deCaName) | Out-Null

$dede.Range("H3").Value = "2016-03-22 00:43:36"

Write-Output "Handback report generated"
